# Append one new row (row 13) of data to Sheet1, mirroring the existing
# rows' "numbers/blank stored as text" convention (ignoredErrors
# numberStoredAsText covers A1:H12 -> A1:H13 after this edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

# A13 is blank (stored as an empty text value, like A2:A4/A6:A12).
# C13 ("222") looks numeric but must stay text, like C2/C3 above it.
# A leading apostrophe forces text storage for both; the trailing
# ".Style = 'Normal'" below then drops the quote-prefix formatting so
# the cells end up with the same (default) style as the rest of the sheet.
$ws.Cells.Item($row, 1).Value = "'"
$ws.Cells.Item($row, 2).Value = "احمد"
$ws.Cells.Item($row, 3).Value = "'222"
$ws.Cells.Item($row, 4).Value = "الصمود"
$ws.Cells.Item($row, 5).Value = "الرحلة 3"
$ws.Cells.Item($row, 6).Value = "C3"
$ws.Cells.Item($row, 7).Value = "NRC"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٥٨:٣٩ م"

$ws.Range("A$row" + ":C$row").Style = "Normal"
